# Applies the commit's changes to the "Blank 3 Statement Model" / "Upload_Guidance" sheets.

$wb = $excel.ActiveWorkbook
$model = $wb.Worksheets.Item("Blank 3 Statement Model")

# --- Relabel assumption rows (A15:A17) ---
$model.Range("A15").Value = "D&A % of Sales"
$model.Range("A16").Value = "LTD Rate (Avg Debt)"
$model.Range("A17").Value = "Tax Rate (Assumption)"

# --- Row 42 (Taxes): replace formulas with static zeros ---
$model.Range("B42:E42").Value = 0

# --- Row 77 (Retained Earnings): replace formulas/value with static zeros ---
$model.Range("B77:E77").Value = 0
